# Fix up the "保險" (insurance) and "債務" (debt) sheets: correct the
# header row (which had been populated with stray data values instead of
# real field names) and fill in the missing metadata columns so both
# sheets match the other property-type sheets (land/building/car/...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 8: 保險 (insurance)
# ---------------------------------------------------------------
$ins = $wb.Worksheets.Item(8)

# Header row (row 1) - fix columns B:D and add E:K
$ins.Range("B1").Value = "company"
$ins.Range("C1").Value = "name"
$ins.Range("D1").Value = "owner"
$ins.Range("E1").Value = "property_category"
$ins.Range("F1").Value = "category"
$ins.Range("G1").Value = "date"
$ins.Range("H1").Value = "legislator_name"
$ins.Range("I1").Value = "legislator_id"
$ins.Range("J1").Value = "source_file"
$ins.Range("K1").Value = "index"

# Row 2
$ins.Range("E2").Value = "insurance"
$ins.Range("F2").Value = "normal"
$ins.Range("G2").NumberFormat = "@"
$ins.Range("G2").Value = "2012-04-27"
$ins.Range("H2").Value = "李昆澤"
$ins.Range("I2").Value = 1327
$ins.Range("J2").Value = "tmp72b91"
$ins.Range("K2").Value = 93

# Row 3
$ins.Range("E3").Value = "insurance"
$ins.Range("F3").Value = "normal"
$ins.Range("G3").NumberFormat = "@"
$ins.Range("G3").Value = "2012-04-27"
$ins.Range("H3").Value = "李昆澤"
$ins.Range("I3").Value = 1327
$ins.Range("J3").Value = "tmp72b91"
$ins.Range("K3").Value = 94

# Row 4
$ins.Range("E4").Value = "insurance"
$ins.Range("F4").Value = "normal"
$ins.Range("G4").NumberFormat = "@"
$ins.Range("G4").Value = "2012-04-27"
$ins.Range("H4").Value = "李昆澤"
$ins.Range("I4").Value = 1327
$ins.Range("J4").Value = "tmp72b91"
$ins.Range("K4").Value = 95

# Match formatting: copy the bold/bordered header style onto the new
# header cells, and the plain data style onto the new data cells.
$ins.Range("B1").Copy()
$ins.Range("E1:K1").PasteSpecial(-4122)
$ins.Range("B2").Copy()
$ins.Range("E2:K2").PasteSpecial(-4122)
$ins.Range("B3").Copy()
$ins.Range("E3:K3").PasteSpecial(-4122)
$ins.Range("B4").Copy()
$ins.Range("E4:K4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Sheet 9: 債務 (debt)
# ---------------------------------------------------------------
$debt = $wb.Worksheets.Item(9)

# Header row (row 1) - fix columns B:G and add H:N
$debt.Range("B1").Value = "species"
$debt.Range("C1").Value = "debtor"
$debt.Range("D1").Value = "owner"
$debt.Range("E1").Value = "total"
$debt.Range("F1").Value = "register_date"
$debt.Range("G1").Value = "register_reason"
$debt.Range("H1").Value = "property_category"
$debt.Range("I1").Value = "category"
$debt.Range("J1").Value = "date"
$debt.Range("K1").Value = "legislator_name"
$debt.Range("L1").Value = "legislator_id"
$debt.Range("M1").Value = "source_file"
$debt.Range("N1").Value = "index"

# Row 2: the existing B2:G2 values were shifted by one column each
# (C/D/E/F/G held the wrong field); put every value back where it
# belongs and append the shared metadata columns H:N.
$debt.Range("B2").Value = "購屋貸款"
$debt.Range("C2").Value = "廖翊百"
$debt.Range("D2").Value = "台灣銀行前金分行高雄市前金區中正四路"
$debt.Range("E2").Value = 4965967
$debt.Range("F2").Value = "100年10月18日"
$debt.Range("G2").Value = "購屋貸款"
$debt.Range("H2").Value = "debt"
$debt.Range("I2").Value = "normal"
$debt.Range("J2").NumberFormat = "@"
$debt.Range("J2").Value = "2012-04-27"
$debt.Range("K2").Value = "李昆澤"
$debt.Range("L2").Value = 1327
$debt.Range("M2").Value = "tmp72b91"
$debt.Range("N2").Value = 105

# Match formatting as above.
$debt.Range("B1").Copy()
$debt.Range("H1:N1").PasteSpecial(-4122)
$debt.Range("B2").Copy()
$debt.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "edit applied"
